$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 450.80768
$ws.Range("I33").Value = 86.09999999999999
$ws.Range("J33").Value = 1666.5
$ws.Range("K33").Value = 86.09999999999999
$ws.Range("L33").Value = 1666.5
$ws.Range("M33").Value = 142.9
$ws.Range("N33").Value = -2124.5
$ws.Range("H43").Value = 1504.6
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1655.75
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1655.75
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1793.75
$ws.Range("H74").Value = 4157.9473
$ws.Range("I74").Value = 4381.909
$ws.Range("J74").Value = 3850
$ws.Range("K74").Value = 4381.909
$ws.Range("L74").Value = 3850
$ws.Range("M74").Value = -3445.909
$ws.Range("N74").Value = -5722
$ws.Range("H77").Value = 4157.9473
$ws.Range("I77").Value = 4381.909
$ws.Range("J77").Value = 3850
$ws.Range("K77").Value = 21909.545
$ws.Range("L77").Value = 19250
$ws.Range("M77").Value = -17229.545
$ws.Range("N77").Value = -28610
$ws.Range("H103").Value = 1217.5
$ws.Range("I103").Value = 833.3333
$ws.Range("K103").Value = 2499.9999
$ws.Range("M103").Value = -1913.9999
$ws.Range("H129").Value = 925.2286
$ws.Range("J129").Value = 1052.3448
$ws.Range("L129").Value = 3157.0344
$ws.Range("N129").Value = -13157.0344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17276934
$ws.Range("I61").Value = 21761828
$ws.Range("J61").Value = 84844
$ws.Range("K61").Value = 21761828
$ws.Range("L61").Value = 84844
$ws.Range("M61").Value = -21761616
$ws.Range("N61").Value = -85268
$ws.Range("H97").Value = 1954239.1
$ws.Range("I97").Value = 2842133.8
$ws.Range("K97").Value = 2842133.8
$ws.Range("M97").Value = -2841637.8
$ws.Range("H132").Value = 49633.07
$ws.Range("I132").Value = 35186.31
$ws.Range("J132").Value = 79558.5
$ws.Range("K132").Value = 105558.93
$ws.Range("L132").Value = 238675.5
$ws.Range("M132").Value = -103028.93
$ws.Range("N132").Value = -243735.5
$ws.Range("H136").Value = 17276934
$ws.Range("I136").Value = 21761828
$ws.Range("J136").Value = 84844
$ws.Range("K136").Value = 65285484
$ws.Range("L136").Value = 254532
$ws.Range("M136").Value = -65282934
$ws.Range("N136").Value = -259632

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 15750
$ws.Range("J81").Value = 15750
$ws.Range("L81").Value = 15750
$ws.Range("N81").Value = -17872
$ws.Range("H84").Value = 15750
$ws.Range("J84").Value = 15750
$ws.Range("L84").Value = 47250
$ws.Range("N84").Value = -57858
$ws.Range("H86").Value = 12361.5
$ws.Range("I86").Value = 15301.883
$ws.Range("J86").Value = 2364.2
$ws.Range("K86").Value = 15301.883
$ws.Range("L86").Value = 2364.2
$ws.Range("M86").Value = -14178.883
$ws.Range("N86").Value = -4610.2
$ws.Range("H89").Value = 12361.5
$ws.Range("I89").Value = 15301.883
$ws.Range("J89").Value = 2364.2
$ws.Range("K89").Value = 76509.41499999999
$ws.Range("L89").Value = 11821
$ws.Range("M89").Value = -70893.41499999999
$ws.Range("N89").Value = -23053
$ws.Range("H99").Value = 931.4286
$ws.Range("I99").Value = 906.125
$ws.Range("J99").Value = 965.1667
$ws.Range("K99").Value = 906.125
$ws.Range("L99").Value = 965.1667
$ws.Range("M99").Value = 591.875
$ws.Range("N99").Value = -3961.1667
$ws.Range("H105").Value = 20002408
$ws.Range("I105").Value = 35716680
$ws.Range("K105").Value = 35716680
$ws.Range("M105").Value = -35714933
$ws.Range("H133").Value = 46235
$ws.Range("J133").Value = 46235
$ws.Range("L133").Value = 46235
$ws.Range("N133").Value = -56355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2806.92
$ws.Range("I31").Value = 1521.7646
$ws.Range("J31").Value = 5537.875
$ws.Range("K31").Value = 1521.7646
$ws.Range("L31").Value = 5537.875
$ws.Range("M31").Value = -1226.7646
$ws.Range("N31").Value = -6127.875
$ws.Range("H34").Value = 2806.92
$ws.Range("I34").Value = 1521.7646
$ws.Range("J34").Value = 5537.875
$ws.Range("K34").Value = 1521.7646
$ws.Range("L34").Value = 5537.875
$ws.Range("M34").Value = -1319.7646
$ws.Range("N34").Value = -5941.875
$ws.Range("H63").Value = 35655.617
$ws.Range("J63").Value = 35655.617
$ws.Range("L63").Value = 35655.617
$ws.Range("N63").Value = -37027.617
$ws.Range("H66").Value = 35655.617
$ws.Range("J66").Value = 35655.617
$ws.Range("L66").Value = 106966.851
$ws.Range("N66").Value = -113830.851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1797.875
$ws.Range("J75").Value = 3015
$ws.Range("L75").Value = 9045
$ws.Range("N75").Value = -11041
$ws.Range("H78").Value = 1797.875
$ws.Range("J78").Value = 3015
$ws.Range("L78").Value = 27135
$ws.Range("N78").Value = -37119
$ws.Range("H113").Value = 477.77777
$ws.Range("I113").Value = 376.66666
$ws.Range("J113").Value = 528.3333
$ws.Range("K113").Value = 1129.99998
$ws.Range("L113").Value = 1584.9999
$ws.Range("M113").Value = 1040.00002
$ws.Range("N113").Value = -5924.9999
$ws.Range("H117").Value = 2899745.2
$ws.Range("I117").Value = 475.8
$ws.Range("J117").Value = 5129952.5
$ws.Range("K117").Value = 1427.4
$ws.Range("L117").Value = 15389857.5
$ws.Range("M117").Value = 2014.6
$ws.Range("N117").Value = -15396741.5
$ws.Range("H131").Value = 1025.2239
$ws.Range("J131").Value = 1121.7241
$ws.Range("L131").Value = 3365.1723
$ws.Range("N131").Value = -13445.1723
$ws.Range("H141").Value = 6848.4546
$ws.Range("I141").Value = 6848.4546
$ws.Range("K141").Value = 20545.3638
$ws.Range("M141").Value = -15365.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10000000
$ws.Range("I20").Value = 10000000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10000000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9999755
$ws.Range("N20").ClearContents()
$ws.Range("H21").Value = 5005000
$ws.Range("I21").Value = 5005000
$ws.Range("K21").Value = 5005000
$ws.Range("M21").Value = -5004827
$ws.Range("H29").Value = 25000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 25000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 25000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -25580
$ws.Range("H30").Value = 5005000
$ws.Range("I30").Value = 5005000
$ws.Range("K30").Value = 5005000
$ws.Range("M30").Value = -5004895
$ws.Range("H70").Value = 64247.53
$ws.Range("I70").Value = 104940.8
$ws.Range("J70").Value = 6114.2856
$ws.Range("K70").Value = 104940.8
$ws.Range("L70").Value = 6114.2856
$ws.Range("M70").Value = -104670.8
$ws.Range("N70").Value = -6654.2856
$ws.Range("H73").Value = 64247.53
$ws.Range("I73").Value = 104940.8
$ws.Range("J73").Value = 6114.2856
$ws.Range("K73").Value = 104940.8
$ws.Range("L73").Value = 6114.2856
$ws.Range("M73").Value = -104004.8
$ws.Range("N73").Value = -7986.2856
$ws.Range("H126").Value = 2589.125
$ws.Range("J126").Value = 2877.1667
$ws.Range("L126").Value = 8631.500100000001
$ws.Range("N126").Value = -13571.5001
$ws.Range("H132").Value = 55660.297
$ws.Range("I132").Value = 38348.85
$ws.Range("J132").Value = 102401.2
$ws.Range("K132").Value = 115046.55
$ws.Range("L132").Value = 307203.6
$ws.Range("M132").Value = -112516.55
$ws.Range("N132").Value = -312263.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6107.077
$ws.Range("I7").Value = 3581.0908
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 3581.0908
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = -3469.0908
$ws.Range("N7").Value = -20224
$ws.Range("H40").Value = 2531.9795
$ws.Range("I40").Value = 2205.878
$ws.Range("J40").Value = 4203.25
$ws.Range("K40").Value = 2205.878
$ws.Range("L40").Value = 4203.25
$ws.Range("M40").Value = -2069.878
$ws.Range("N40").Value = -4475.25
$ws.Range("H48").Value = 13001
$ws.Range("J48").Value = 13001
$ws.Range("L48").Value = 13001
$ws.Range("N48").Value = -14323
$ws.Range("H126").Value = 6107.077
$ws.Range("I126").Value = 3581.0908
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 10743.2724
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -8273.2724
$ws.Range("N126").Value = -64940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3668.1667
$ws.Range("I62").Value = 4333.3335
$ws.Range("J62").Value = 3003
$ws.Range("K62").Value = 4333.3335
$ws.Range("L62").Value = 3003
$ws.Range("M62").Value = -3709.3335
$ws.Range("N62").Value = -4251
$ws.Range("H65").Value = 3668.1667
$ws.Range("I65").Value = 4333.3335
$ws.Range("J65").Value = 3003
$ws.Range("K65").Value = 21666.6675
$ws.Range("L65").Value = 15015
$ws.Range("M65").Value = -18546.6675
$ws.Range("N65").Value = -21255
$ws.Range("H107").Value = 397.73685
$ws.Range("I107").Value = 355.2857
$ws.Range("J107").Value = 516.6
$ws.Range("K107").Value = 1065.8571
$ws.Range("L107").Value = 1549.8
$ws.Range("M107").Value = 854.1428999999998
$ws.Range("N107").Value = -5389.8
